# Weekly data refresh: a new price observation is inserted as the new
# row 166 (pushing the existing rows 166-246 down to 167-247). The new
# row starts out as a duplicate of the row that lands at 167 (the old
# row 166) and then gets its own date and origin.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(166).Insert()

$ws.Range("A167:R167").Copy()
$ws.Range("A166").PasteSpecial()

$ws.Range("D166").Value = 44755
$ws.Range("O166").Value = "Región de Arica y Parinacota"
